# Auto-generated PowerShell COM-interop script
# Applies: status (STATUS_ABERTURA) updates on existing rows,
# a couple of cell backfills, and 62 new production rows (520-581).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($rng, $text)
    # Force pure-numeric-looking strings (e.g. CNPJ) to be stored as text,
    # matching the workbook's existing convention for column B, while
    # resetting the style afterwards so no stray number-format/quote-prefix
    # style sticks to the cell.
    $rng.Value = "'" + $text
    $rng.Style = 'Normal'
}

# --- STATUS_ABERTURA (column F) updates on existing rows ---
$statusUpdates = @{
    19 = 'APROVADA'
    77 = 'REPROVADA'
    132 = 'PENDÊNCIA DOC'
    160 = 'REPROVADA'
    241 = 'APROVADA'
    325 = 'APROVADA'
    453 = 'Ainda nao iniciou a abertura de conta'
    455 = 'Ainda nao iniciou a abertura de conta'
    456 = 'Ainda nao iniciou a abertura de conta'
    457 = 'Ainda nao iniciou a abertura de conta'
    469 = 'Ainda nao iniciou a abertura de conta'
    470 = 'Ainda nao iniciou a abertura de conta'
    472 = 'Ainda nao iniciou a abertura de conta'
    473 = 'Ainda nao iniciou a abertura de conta'
    474 = 'Ainda nao iniciou a abertura de conta'
    485 = 'Ainda nao iniciou a abertura de conta'
    489 = 'Ainda nao iniciou a abertura de conta'
    493 = 'APROVADA'
    499 = 'Ainda nao iniciou a abertura de conta'
    502 = 'Ainda nao iniciou a abertura de conta'
    509 = 'Ainda nao iniciou a abertura de conta'
    511 = 'Ainda nao iniciou a abertura de conta'
    513 = 'Ainda nao iniciou a abertura de conta'
    514 = 'Ainda nao iniciou a abertura de conta'
}
foreach ($r in $statusUpdates.Keys) {
    $ws.Range("F$r").Value = $statusUpdates[$r]
}

# --- Misc single-cell backfills ---
# Row 241: PENDENCIAS goes back to "no pendencies" (0) now that the record is approved
$ws.Range("G241").Value = 0
# Row 441: PENDENCIAS was left blank; set it to 0
$ws.Range("G441").Value = 0
# Row 517: CONSULTOR / ORIGEM were left blank; backfill them
$ws.Range("D517").Value = "Caua Miedes do Espirito Santo"
$ws.Range("E517").Value = "URA"

# --- New production rows (520-581) ---
$newRows = @(
    @{Row=520; Date=46063; B='63438186000185'; C='J&V STORE LTDA'; D='Julia do Nascimento Lopes'; E='URA'; F='APROVADA'; G=0}
    @{Row=521; Date=46063; B='54665244000112'; C='TRANSPORTE CRISTINA LTDA'; D='Maria Luisa Ribeiro da Silva'; E='URA'; F='ANÁLISE'; G=0}
    @{Row=522; Date=46063; B='11233505000112'; C='ROSY ESTEVES ENCARNACAO LTDA'; D='Joao Pedro Gabriel Troiano'; E='REPESCAGEM'; F='APROVADA'; G=0}
    @{Row=523; Date=46063; B='22736165000169'; C='JEAN CARLOS AVELINO PEREIRA DE BARROS'; D='Ana Beatriz Rodrigues'; E='URA'; F='ANÁLISE'; G=0}
    @{Row=524; Date=46063; B='39145774000146'; C='QUALITY CONSTRUCOES E SERVICOS LTDA'; D='Jhonatan Vinicius de Paula Alonso'; E='URA'; F='ANÁLISE'; G=0}
    @{Row=525; Date=46063; B='14779669000138'; C='AMAZON GEOGRAPHIC HOTEIS E TURISMO LTDA'; D='Felipe Arantes Martelo'; E='REPESCAGEM'; F='ANÁLISE'; G=0}
    @{Row=526; Date=46063; B='25333086000103'; C='MARIA APARECIDA PADRE DA SILVA SOUSA 08611199405'; D='Mariana Gabriela Ferreira Barboza'; E='REPESCAGEM'; F='APROVADA'; G=0}
    @{Row=527; Date=46063; B='14972560000112'; C='KM JEANS CONFECCOES LTDA'; D='Miriam Helena Franzoi'; E='REPESCAGEM'; F='ANÁLISE'; G=0}
    @{Row=528; Date=46063; B='33988761000152'; C='FILIPE ROSA DA SILVA ENERGIA SOLAR'; D='Otavio Henrique Silva de Avila'; E='REPESCAGEM'; F='APROVADA'; G=0}
    @{Row=529; Date=46063; B='44131714000104'; C='PARK INTERATIVO EVENTOS LTDA'; D='Gleicy Ferreira Gomes'; E='REPESCAGEM'; F='APROVADA'; G=0}
    @{Row=530; Date=46063; B='64972929000165'; C='BRASIL & USA TRANSPORTES LTDA'; D='Bruna Tobias Dos Santos'; E='DISCADOR'; F='APROVADA'; G=0}
    @{Row=531; Date=46063; B='50698537000145'; C='SERVICO DE INSPECAO VEICULAR KM 32 LTDA'; D='Miriam Helena Franzoi'; E='REPESCAGEM'; F='ANÁLISE'; G=0}
    @{Row=532; Date=46063; B='64854886000113'; C='GLL SOLUCOES CONTABEIS LTDA'; D='Jhonatan Vinicius de Paula Alonso'; E='REPESCAGEM'; F='ANÁLISE'; G=0}
    @{Row=533; Date=46063; B='57801824000188'; C='ISA RANGEL COMERCIO VAREJISTA DE ROUPAS LTDA'; D='Giovana Vitoria da Silva'; E='URA'; F='PENDÊNCIA DOC'; G='Procuracao com poderes de abrir e movimentar conta corrente junto a instituicoes financeiras com mandato vigente e devidamente assinada<br><br>Contrato Social atualizado e registrado no orgao competente'}
    @{Row=534; Date=46063; B='58419682000151'; C='MASTER EM CONSORCIOS LTDA'; D='Marcelo Sanches Espirito Santo'; E='URA'; F='APROVADA'; G=0}
    @{Row=535; Date=46063; B='64813988000190'; C='GMDC CONSULTORIA AMBIENTAL LTDA'; D='Endrew Lourenco Posca'; E='DISCADOR'; F='ANÁLISE'; G=0}
    @{Row=536; Date=46063; B='15815893000109'; C='MARCENARIA AMAZONAS LIMITADA'; D='Giovana Vitoria da Silva'; E='URA'; F='ANÁLISE'; G=$null}
    @{Row=537; Date=46063; B='09308808000105'; C='DIAS E LIRA LTDA'; D='Ana Beatriz Rodrigues'; E='URA'; F='APROVADA'; G=0}
    @{Row=538; Date=46063; B='28288285000127'; C='CONSTRUVILLE ENGENHARIA E CONSTRUCOES LTDA'; D='Franciane Roberta Cordeiro'; E='URA'; F='REPROVADA'; G=0}
    @{Row=539; Date=46063; B='64923951000115'; C='CARRO CERTO VEICULOS LTDA'; D='Cristiane dos Santos Andrade'; E='LEAD MANUAL'; F='APROVADA'; G=0}
    @{Row=540; Date=46063; B='64904270000100'; C='DANIEL BARBOSA TERRA LTDA'; D='Ana Carla Ferreira Fellippini'; E='DISCADOR'; F='APROVADA'; G=0}
    @{Row=541; Date=46063; B='38949358000138'; C='FORMULA USINAGEM LTDA'; D='Franciane Roberta Cordeiro'; E='URA'; F='ANÁLISE'; G=0}
    @{Row=542; Date=46063; B='43689760000160'; C='POWER TECH GMC LTDA'; D='Caua Miedes do Espirito Santo'; E='URA'; F='ANÁLISE'; G=0}
    @{Row=543; Date=46063; B='60748069000120'; C='EDY FABRICACAO E MONTAGENS LTDA'; D='Franciane Roberta Cordeiro'; E='URA'; F='APROVADA'; G=0}
    @{Row=544; Date=46063; B='58970777000169'; C='WELISON CASTIONI FERREIRA'; D='Ana Clara Sabio de Souza'; E='REPESCAGEM'; F='APROVADA'; G=0}
    @{Row=545; Date=46063; B='31245109000112'; C='TAISA MIRELA MESSIAS DE OLIVEIRA'; D='Ana Beatriz Rodrigues'; E='URA'; F='REPROVADA'; G=0}
    @{Row=546; Date=46063; B='40129863000188'; C='GILSTON DE PAULA'; D='Giovana Vitoria da Silva'; E='URA'; F='APROVADA'; G=0}
    @{Row=547; Date=46063; B='64753679000172'; C='TATIANA FREIRE SERVICOS ADMINISTRATIVOS'; D='Maria Luisa Ribeiro da Silva'; E='URA'; F='APROVADA'; G=0}
    @{Row=548; Date=46063; B='64863917000100'; C='NALDO REPRESENTACOES LTDA'; D='Ana Laura Rodrigues da Silva'; E='DISCADOR'; F='APROVADA'; G=0}
    @{Row=549; Date=46063; B='52102556000192'; C='MONTA RJ - MONTAGENS INDUSTRIAIS LTDA'; D='Felipe Arantes Martelo'; E='URA'; F='ANÁLISE'; G=0}
    @{Row=550; Date=46063; B='64980093000140'; C='BOM - SERVICOS MEDICOS DE APOIO A SAUDE LTDA'; D='Endrew Lourenco Posca'; E='DISCADOR'; F='ANÁLISE'; G=0}
    @{Row=551; Date=46063; B='47142388000129'; C='MIRIAN BARBOSA DA SILVA'; D='Felipe Arantes Martelo'; E='URA'; F='APROVADA'; G=0}
    @{Row=552; Date=46063; B='08491873000156'; C='TAKAMORI GESTAO EMPRESARIAL LTDA'; D='Jhonatan Vinicius de Paula Alonso'; E='REPESCAGEM'; F='APROVADA'; G=0}
    @{Row=553; Date=46063; B='49629584000158'; C='MAAF CONSULTORIA LTDA'; D='Bruna Tobias Dos Santos'; E='URA'; F='APROVADA'; G=0}
    @{Row=554; Date=46063; B='39604808000113'; C='PRESENZA CAR AUTO CENTER LTDA'; D='Sofia Helena Vieira Domingues'; E='URA'; F='APROVADA'; G=0}
    @{Row=555; Date=46063; B='46545911000103'; C='POUSADA BRASIL TURISMO LTDA'; D='Marcelo Sanches Espirito Santo'; E='URA'; F='APROVADA'; G=0}
    @{Row=556; Date=46063; B='26337655000143'; C='JOSE ALFREDO MUSSI DE SOUZA & CIA LTDA'; D='Ana Beatriz Rodrigues'; E='TPG'; F='REPROVADA'; G=0}
    @{Row=557; Date=46063; B='18222665000105'; C='MADEREIRA ALIANCA LTDA'; D='Franciane Roberta Cordeiro'; E='URA'; F='APROVADA'; G=0}
    @{Row=558; Date=46063; B='13461347000183'; C='LUAU PRODUCOES E EVENTOS LIMITADA'; D='Joao Pedro Gabriel Troiano'; E='REPESCAGEM'; F='ANÁLISE'; G=0}
    @{Row=559; Date=46063; B='57505002000150'; C='57.505.002 JOSELI VIEIRA DOS SANTOS SILVA'; D='Ana Beatriz Rodrigues'; E='TPG'; F='APROVADA'; G=0}
    @{Row=560; Date=46063; B='65000115000121'; C='B&R IMOBILIARIOS LTDA'; D='Endrew Lourenco Posca'; E='DISCADOR'; F='ANÁLISE'; G=0}
    @{Row=561; Date=46063; B='40095178000188'; C='BARBER HAUSS BARBEARIA LTDA'; D='Tauani Santos de Andrade'; E='URA'; F='APROVADA'; G=0}
    @{Row=562; Date=46063; B='64289840000107'; C='NP REPRESENTACOES COMERCIAL LTDA'; D='Yara Galicia de Andrade dos Santos'; E='URA'; F='APROVADA'; G=0}
    @{Row=563; Date=46063; B='52889231000100'; C='SS ROCHA ASSESSORIA E CONSULTORIA LTDA'; D='Jamily de Lima Alves dos Santos'; E='URA'; F='APROVADA'; G=0}
    @{Row=564; Date=46063; B='64986006000162'; C='POUSADA E KIOSQUE DAG E JULIO LTDA'; D='Bruna Tobias Dos Santos'; E='DISCADOR'; F='ANÁLISE'; G=0}
    @{Row=565; Date=46063; B='60813940000121'; C='ITAVERAVA TRANSPORTES'; D='Yara Galicia de Andrade dos Santos'; E='URA'; F='APROVADA'; G=0}
    @{Row=566; Date=46063; B='41414527000140'; C='VANESSA BARBOSA CHIOVETTI MORAES CONSULTORIO ODONTOLOGICO LTDA'; D='Maria Luisa Ribeiro da Silva'; E='URA'; F='PENDÊNCIA DOC'; G=0}
    @{Row=567; Date=46063; B='39964872000105'; C='39.964.872 NATALIA FAUSTO DINIZ MARTINS'; D='Ana Beatriz Rodrigues'; E='TPG'; F='APROVADA'; G=0}
    @{Row=568; Date=46063; B='62360809000181'; C='INSTITUTO TEMPO DE CUIDAR CARETIME LTDA'; D='Caua Miedes do Espirito Santo'; E='URA'; F='APROVADA'; G=0}
    @{Row=569; Date=46063; B='62411121000183'; C='K&G AUTOMACAO INDUSTRIAL LTDA'; D='Sofia Helena Vieira Domingues'; E='URA'; F='REPROVADA'; G=0}
    @{Row=570; Date=46063; B='64805962000109'; C='MV VIDROS LTDA'; D='Endrew Lourenco Posca'; E='DISCADOR'; F='REPROVADA'; G=0}
    @{Row=571; Date=46063; B='11812991000123'; C='MAQJOB LOCACAO DE MAQUINAS LTDA'; D='Miriam Helena Franzoi'; E='URA'; F='PENDÊNCIA DOC'; G='Procuracao com poderes de abrir e movimentar conta corrente junto a instituicoes financeiras com mandato vigente e devidamente assinada<br><br>Contrato Social atualizado e registrado no orgao competente'}
    @{Row=572; Date=46063; B='60612366000143'; C='EMINENTES FORMATURAS LTDA'; D='Maria Luisa Ribeiro da Silva'; E='URA'; F='REPROVADA'; G=0}
    @{Row=573; Date=46063; B='34667260000137'; C='LOPES BAR E LANCHONETE LTDA'; D='Caua Miedes do Espirito Santo'; E='URA'; F='APROVADA'; G=0}
    @{Row=574; Date=46063; B='20667781000180'; C='ANDREIA MARIA PIMENTEL'; D='Giovana Vitoria da Silva'; E='URA'; F='APROVADA'; G=0}
    @{Row=575; Date=46063; B='37129220000120'; C='O S J CONSTRUCOES LTDA'; D='Maria Luisa Ribeiro da Silva'; E='REPESCAGEM'; F='ANÁLISE'; G=0}
    @{Row=576; Date=46064; B='64990822000140'; C='JARDEL SANTANA LTDA'; D='Tauani Santos de Andrade'; E='DISCADOR'; F='APROVADA'; G=0}
    @{Row=577; Date=46064; B='15811504000169'; C='TATIANE KELEN TEBAR'; D='Ana Beatriz Rodrigues'; E='TPG'; F='APROVADA'; G=0}
    @{Row=578; Date=46064; B='44961087000139'; C='GRUPO NOVA ERA COMERCIO DE AUTOMOVEIS LTDA'; D='Miriam Helena Franzoi'; E='URA'; F='APROVADA'; G=0}
    @{Row=579; Date=46064; B='63835109000169'; C='E M T COMERCIO DE MAQUINAS E SERVICE AMBIENTAIS LTDA'; D='Stephany Eduarda Pereira'; E='DISCADOR'; F='APROVADA'; G=0}
    @{Row=580; Date=46064; B='35196665000105'; C='THIAGO JUNIO ALVES DA SILVA '; D='Miriam Helena Franzoi'; E='URA'; F='APROVADA'; G=0}
    @{Row=581; Date=46064; B='50217786000171'; C='PRIME SERVICE ASSEIO E CONSERVACAO LTDA'; D='Jamily de Lima Alves dos Santos'; E='URA'; F='APROVADA'; G=0}
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.Date
    Set-TextValue $ws.Range("B$r") $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
    if ($null -ne $item.G) {
        $ws.Range("G$r").Value = $item.G
    }
}

# --- Restore the active-cell selection left by the author ---
$ws.Range("C11").Select()
